$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'lunges pad',
    'knee protection pad',
    'boys basketball pants',
    'working knee pads for men',
    'knee yoga pants',
    'padded knee sleeve for sliding',
    'men tight pants',
    'protective compression wear',
    'spandex compression shorts men',
    'softball mens pants',
    'compression capri',
    'softball sliding pants youth girls',
    'baseball pants black mens',
    'wrestling knee pads pair',
    'compressions tights for men',
    'men sports leggings',
    'compression shorts boys padded',
    'basketball padded knee sleeve',
    'yoga pad thick',
    'mens leggings tall',
    'knee protector construction',
    'basketball leggings for girls',
    'lacrosse compression shorts',
    'boys sports tights leggings',
    'mens sweat pads',
    'mens construction knee pads',
    'yoga pants men',
    'bursitis knee',
    'thread protector paintball',
    'tight gym pants men',
    'adult bee tights',
    'basketball shorts in bulk',
    'rash guard men leggings',
    'youth sliding short',
    'running knee compression',
    'knee protection for running',
    'girls basketball leggings',
    'youth xxl baseball pants',
    'boys softball pants',
    'tight for boys',
    'big and tall leggings men',
    'knee protect',
    'knee sleeve basketball men',
    'protective baseball',
    'compression pants for recovery',
    'girls basketball knee guards',
    'basketball knee pads for women',
    'gel wrestling knee pads',
    'mens hockey pads',
    'long knee pads volleyball',
    'mens compression pants 3/4 length',
    'knee pads motorcycle',
    'sliding shorts youth girls',
    'boys athletic leggings youth',
    'knee pads for basketball women',
    'big mens compression pants',
    'soccer sliding pants',
    'good thread mens pants',
    'xl knee pads for men',
    'capri pouches adults',
    'knee sleeves hex',
    'knee compression sleeve with padding',
    'baleaf mens pants',
    'leggings tight',
    'compression tights mens',
    'compression knee sleeve basketball',
    'volleyball knee pads gel',
    'knee pads for work xxl',
    'knee pad for exercise',
    'running tights youth boys',
    'boy capri pants',
    'youth hockey girdle',
    'youth compression leggings boys',
    'sliding compression shorts',
    'baseball leg protection',
    'knee sleeves for basketball youth',
    'youth paintball pants',
    'youth boy tights',
    'softball catcher pants',
    'knee pads volleyball mens',
    'mens running compression pants',
    'football tights youth boys',
    'athletic tights youth boys',
    'padded calf sleeve',
    'mens softball pants',
    'mens baseball shorts',
    'knee pad for soccer',
    'compression leggings boys',
    'capri for men',
    'mens workout pants leggings',
    'knee pad sport',
    'tights pants boys',
    'softball sliding shorts women',
    'padded sliding shorts',
    'arthritis hope knee sleeve',
    'padded knee sleeve basketball',
    'snowboarding mens pants',
    'knee compression cold',
    'knee pads squats',
    'men work pants with knee pads'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
